# Add validation text for "create slot" rules:
#   - "Slot must be create at least 6 hours from now" ->
#       append " and maximum 30 days from now"
#   - "Each slot must be at least 15 minutes long" ->
#       append " and maximum 2 hours 15 minutes long"
#
# Both additions must land in their own new <w:r> (matching the target
# OOXML diff) with the same run formatting as the run they follow
# (Times New Roman / sz 22 / szCs 22 / lang en-US). A plain
# Range.InsertAfter() call would just extend the existing run (and
# save-time canonicalization re-merges runs with identical formatting),
# so the new text is inserted via Range.InsertXML with an explicit
# <w:r> fragment -- that is the one path that reliably produces a
# distinct sibling run even when its rPr is identical to its neighbor.

function Add-SentenceToParagraph {
    param(
        [string]$containsText,
        [string]$appendText
    )

    $d = $word.ActiveDocument

    $target = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like ("*" + $containsText + "*")) {
            $target = $p
            break
        }
    }

    if ($target -eq $null) {
        throw "Paragraph containing '$containsText' was not found"
    }

    # Position right before the paragraph mark (End is one past the
    # last character, the last slot being the paragraph mark itself).
    $insertPos = $target.Range.End - 1
    $collapsed = $d.Range($insertPos, $insertPos)

    # Plain insert first so Start/End math for the new run is trivial.
    $collapsed.InsertAfter($appendText)

    $newRunRange = $d.Range($insertPos, $insertPos + $appendText.Length)

    $openXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">$appendText</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

    # InsertXML replaces the (exact, already-inserted-text) range's
    # contents with the parsed run -- giving it its own independent
    # <w:r>/<w:rPr> instead of being coalesced into the previous run.
    $newRunRange.InsertXML($openXml)
}

Add-SentenceToParagraph "Slot must be create at least 6 hours from now" " and maximum 30 days from now"
Add-SentenceToParagraph "Each slot must be at least 15 minutes long" " and maximum 2 hours 15 minutes long"

Write-Output "done"
